$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the capitalisation of the isSourceOf / isDerivedFrom headings
$ws.Range("H1").Value = "IsSourceOf"
$ws.Range("I1").Value = "IsDerivedFrom"

# Update the active selection to I2
$ws.Range("I2").Select()
